# Update the public EPEX Spot prices workbook with the latest day of data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column AP (25-jul) with 24 hourly prices
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# New header cell, copying the formatting of the previous header (AO1)
$wsPrix.Range("AP1").Value = "25-jul"
$wsPrix.Range("AO1").Copy() | Out-Null
$wsPrix.Range("AP1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$prixValues = @{
    2  = 88.69
    3  = 96.52
    4  = 78.48
    5  = 58.99
    6  = 51.97
    7  = 63.43
    8  = 77.97
    9  = 85.98
    10 = 98.34999999999999
    11 = 86.44
    12 = 62.4
    13 = 65.73999999999999
    14 = 59.01
    15 = 35.04
    16 = 16.2
    17 = 40.81
    18 = 46.49
    19 = 58.52
    20 = 63.23
    21 = 94.90000000000001
    22 = 85
    23 = 104.72
    24 = 105
    25 = 81.95999999999999
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Range("AP$row").Value = $prixValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append row 39 with the new daily price
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A39").NumberFormat = "@"
$wsGaz.Range("A39").Value = "2025-07-23"
$wsGaz.Range("A39").ClearFormats()
$wsGaz.Range("B39").Value = 32.15

# ---------------------------------------------------------------------
# Sheet "CO2": append row 39 with the new daily price
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A39").NumberFormat = "@"
$wsCo2.Range("A39").Value = "2025-07-23"
$wsCo2.Range("A39").ClearFormats()
$wsCo2.Range("B39").Value = 68.40000000000001

Write-Host "Workbook updated"
